# Updates cryptos list price (D) and volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.753.15"
$ws.Range("E2").Value = "  +2.45%  "

$ws.Range("D3").Value = "3.742.18"
$ws.Range("E3").Value = "  +6.71%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "419.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.42%  "

$ws.Range("D7").Value = "3.733.09"
$ws.Range("E7").Value = "  +6.74%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.652"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.769"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("E11").Value = "  +11.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000402"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +53.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.75"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.04%  "

$ws.Range("D15").Value = "4.321.22"
$ws.Range("E15").Value = "  +6.55%  "

$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.33%  "

$ws.Range("D18").Value = "3.745.16"
$ws.Range("E18").Value = "  +6.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.10%  "

$ws.Range("E20").Value = "  +3.82%  "

$ws.Range("D21").Value = "66.859.21"
$ws.Range("E21").Value = "  +2.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "444.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +24.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("E25").Value = "  -1.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.29%  "

$ws.Range("E27").Value = "  +2.96%  "

$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.85%  "

$ws.Range("E30").Value = "  +8.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.51%  "

$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D39").Value = "0.0₃0745"
$ws.Range("E39").Value = "  +4.55%  "

$ws.Range("E40").Value = "  -1.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +29.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "29.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +35.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +32.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "146.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "

$ws.Range("E48").Value = "  -4.17%  "

$ws.Range("E49").Value = "  -5.96%  "

$ws.Range("E50").Value = "  -4.50%  "

$ws.Range("E51").Value = "  -2.37%  "
